$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization keys to show grid info (origin/axes/quadrants) on pointer.
# Entered in the same order as authored: keys + most values first, then the
# axis labels (X - Axis / Y - Axis) filled in afterwards.

$ws.Cells.Item(15, 1).Value = "origin"
$ws.Cells.Item(15, 2).Value = "Origin"
$ws.Cells.Item(15, 2).WrapText = $true

$ws.Cells.Item(16, 1).Value = "axisX"
$ws.Cells.Item(17, 1).Value = "axisY"

$ws.Cells.Item(18, 1).Value = "quadrant1"
$ws.Cells.Item(18, 2).Value = "Quadrant 1"
$ws.Cells.Item(18, 2).WrapText = $true

$ws.Cells.Item(19, 1).Value = "quadrant2"
$ws.Cells.Item(19, 2).Value = "Quadrant 2"
$ws.Cells.Item(19, 2).WrapText = $true

$ws.Cells.Item(20, 1).Value = "quadrant3"
$ws.Cells.Item(20, 2).Value = "Quadrant 3"
$ws.Cells.Item(20, 2).WrapText = $true

$ws.Cells.Item(21, 1).Value = "quadrant4"
$ws.Cells.Item(21, 2).Value = "Quadrant 4"
$ws.Cells.Item(21, 2).WrapText = $true

# Fill in the axis labels last (appends new shared strings at the end)
$ws.Cells.Item(16, 2).Value = "X - Axis"
$ws.Cells.Item(16, 2).WrapText = $true

$ws.Cells.Item(17, 2).Value = "Y - Axis"
$ws.Cells.Item(17, 2).WrapText = $true

# Update the active selection to match the final edited cell
$ws.Range("B18").Select()
